$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.009.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.65%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.639.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.06%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.20%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'214.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.44%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5060"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.27%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.009"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.22%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.06443"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.13%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.2573"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.11%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.42%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07706"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.92%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.645.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.82%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.242"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.11%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.864.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.5445"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.62%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₅7901"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.95%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'63.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.37%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'25.999.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.83%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.009"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.25%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'203.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'4.283"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.57%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'9.973"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.03%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.41%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.28%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.941"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +10.29%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'141.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.57%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.1151"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.21%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'15.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.11%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'6.715"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.70%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.05051"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.24%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.241"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.15%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.54%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.190"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.20%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.538"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.49%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.344"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.89%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.638"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -4.63%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.8871"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.12%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = "'1.147.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.83%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = "'0.5609"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.90%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01570"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.37%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.565"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.38%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.009"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.22%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'5.665"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.04%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.8085"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.40%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'99.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.776.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.05%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0₈113"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.83%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.4528"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'1.008"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.09%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'54.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.99%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.05036"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.97%  "
$ws.Range("E51").Style = "Normal"
